# Update odds values on Sheet1 to reflect the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 (ANTIGUA & BARBUDA - ABFA PREMIER LEAGUE: Five Islands - Old Road)
$ws.Range("G5").Value = 22
$ws.Range("H5").Value = 9.5
$ws.Range("J5").Value = 14
$ws.Range("K5").Value = 4.2
$ws.Range("L5").Value = 1.22
$ws.Range("P5").Value = 12
$ws.Range("S5").Value = 1.05
$ws.Range("T5").Value = 9
$ws.Range("U5").Value = 1.6
$ws.Range("V5").Value = 2.26
$ws.Range("W5").Value = 175
$ws.Range("X5").Value = 500
$ws.Range("Y5").Value = 110
$ws.Range("AA5").Value = 400
$ws.Range("AE5").Value = 40
$ws.Range("AF5").Value = 90
$ws.Range("AH5").Value = 27
$ws.Range("AJ5").Value = 16.5
$ws.Range("AL5").Value = 11.25
$ws.Range("AM5").Value = 26
$ws.Range("AN5").Value = 27
$ws.Range("AO5").Value = 150
$ws.Range("AP5").Value = 70
$ws.Range("AS5").Value = 450
$ws.Range("AT5").Value = 8
$ws.Range("AU5").Value = 11.5
$ws.Range("AV5").Value = 55
$ws.Range("AX5").Value = 4.35
$ws.Range("AZ5").Value = 11.25
$ws.Range("BA5").Value = 6.1
$ws.Range("BC5").Value = 80

# Row 8 (ARGENTINA - TORNEO BETANO: Belgrano - Ind. Rivadavia)
$ws.Range("G8").Value = 1.9
$ws.Range("I8").Value = 4.5
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("X8").Value = 7.5
$ws.Range("AH8").Value = 9.5
$ws.Range("BD8").Value = 126

# Row 9 (ARGENTINA - PRIMERA NACIONAL: Gimnasia Mendoza - San Martin T.)
$ws.Range("H9").Value = 2.75
$ws.Range("I9").Value = 3.1
$ws.Range("S9").Value = 1.75
$ws.Range("T9").Value = 2.05
$ws.Range("AE9").Value = 23
$ws.Range("AU9").Value = 11
$ws.Range("BA9").Value = 81
$ws.Range("BB9").Value = 151

# Row 15 (BOLIVIA - DIVISION PROFESIONAL: Oriente Petrolero - Santa Cruz)
$ws.Range("M15").Value = 1.04
$ws.Range("N15").Value = 13
